$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "44.013.89"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "2.242.24"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "267.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +12.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.614"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.87"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.60"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.57%  "
$ws.Range("E13").Value = "  +2.84%  "
$ws.Range("D14").Value = "2.581.92"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.04"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "2.223.73"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.799"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "44.013.51"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("E22").Value = "  +4.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.95"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.54"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +15.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.39"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0919"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.74"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("E36").Value = "  +4.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0356"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +17.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.12"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("E43").Value = "  +1.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0989"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.39"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.50"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.73%  "
$ws.Range("E48").Value = "  +5.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.13"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.441"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.91%  "
